$p = $ppt.ActivePresentation
$nd = $p.Designs.Add()
Write-Host "count $($p.Designs.Count)"
